$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new text values in the precise order they were first authored ---
# (this controls the order new entries are appended to the shared-strings table)
$ws.Range("G6").Value = "3 question daily"
$ws.Range("G9").Value = "yes"
$ws.Range("H7").Value = "Date"
$ws.Range("B3").Value = "Solve the problem in chuck form divide into 2 part"
$ws.Range("F19").Value = "good"

# --- Fill in the remaining "yes" cells in column G (reuse the string) ---
$ws.Range("G10").Value = "yes"
$ws.Range("G11").Value = "yes"
$ws.Range("G12").Value = "yes"
$ws.Range("G13").Value = "yes"
$ws.Range("G14").Value = "yes"
$ws.Range("G16").Value = "yes"
$ws.Range("G19").Value = "yes"
$ws.Range("G21").Value = "yes"

# --- G7: "Revision" (existing shared string) ---
$ws.Range("G7").Value = "Revision"

# --- Apply yellow highlight fill to the "yes" cells (and G17, which gets the fill but no value) ---
$ws.Range("G9").Interior.Color = 65535
$ws.Range("G10").Interior.Color = 65535
$ws.Range("G11").Interior.Color = 65535
$ws.Range("G12").Interior.Color = 65535
$ws.Range("G13").Interior.Color = 65535
$ws.Range("G14").Interior.Color = 65535
$ws.Range("G16").Interior.Color = 65535
$ws.Range("G17").Interior.Color = 65535
$ws.Range("G19").Interior.Color = 65535

# --- Highlight F17 (already bold) with the same yellow fill ---
$ws.Range("F17").Interior.Color = 65535

# --- Fill in "easy" in column E where it was previously blank ---
$ws.Range("E10").Value = "easy"
$ws.Range("E11").Value = "easy"
$ws.Range("E15").Value = "easy"
$ws.Range("E16").Value = "easy"

# --- New Date column (H), formatted like the existing date cells (format code 14) ---
$ws.Range("H9").Value = "45071"
$ws.Range("H10").Value = "45071"
$ws.Range("H11").Value = "45071"
$ws.Range("H12").Value = "45072"
$ws.Range("H13").Value = "45072"
$ws.Range("H14").Value = "45073"
$ws.Range("H16").Value = "45073"
$ws.Range("H17").Value = "45074"
$ws.Range("H19").Value = "45074"
$ws.Range("H21").Value = "45075"
$ws.Range("H22").Value = "45077"
$ws.Range("H25").Value = "45049"

$dateCells = "H9,H10,H11,H12,H13,H14,H16,H17,H19,H21,H22,H25"
foreach ($addr in $dateCells.Split(",")) {
    $ws.Range($addr).NumberFormat = "m/d/yyyy"
}

# --- Column width adjustments ---
$ws.Columns.Item(1).ColumnWidth = 9.5546875
$ws.Columns.Item(2).ColumnWidth = 86.44140625
$ws.Columns.Item(4).ColumnWidth = 13.77734375
$ws.Columns.Item(6).ColumnWidth = 5.33203125
$ws.Columns.Item(8).ColumnWidth = 10.77734375

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("B25").Select()
